$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, shifting existing rows 193:227 down to 194:228
$ws.Rows.Item(193).EntireRow.Insert()

# Populate the newly inserted row 193 with the new weekly data point
$ws.Cells.Item(193, 1).Value = 10
$ws.Cells.Item(193, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value = "La Araucanía"
$ws.Cells.Item(193, 4).Value = 44637
$ws.Cells.Item(193, 5).Value = 9
$ws.Cells.Item(193, 6).Value = 100112039
$ws.Cells.Item(193, 7).Value = "Ciboulette"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 40
$ws.Cells.Item(193, 11).Value = 5000
$ws.Cells.Item(193, 12).Value = 5000
$ws.Cells.Item(193, 13).Value = 5000
$ws.Cells.Item(193, 14).Value = "`$/docena de atados"
$ws.Cells.Item(193, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(193, 16).Value = 1667
$ws.Cells.Item(193, 17).Value = 3
$ws.Cells.Item(193, 18).Value = "Hortaliza"
